$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update location labels to "City, Country" formatting
# (order matches the order new shared strings were introduced)
$ws.Range("B12").Value = "Wuhan, China"
$ws.Range("B7").Value = "Geneva, Switzerland"
$ws.Range("B15").Value = "Zurich, Switzerland"

# Row 17 was LA_CA1 ageband; replace with NYS1 ageband data
$ws.Range("A17").Value = "NYS1"
$ws.Range("B17").Value = "New York State, USA"
$ws.Range("C17").Value = "ageband"
$ws.Range("D17").Value = "data/derived/USA/NYS1_agebands.RDS"
$ws.Range("E17").Value = "yes"

# Row 18 was NYS1 ageband; turn it into NYS1 region data
$ws.Range("A18").Value = "NYS1"
$ws.Range("B18").Value = "New York State, USA"
$ws.Range("C18").Value = "region"
$ws.Range("D18").Value = "data/derived/USA/NYS1_regions.RDS"
$ws.Range("E18").Value = "yes"

# Row 19 was a duplicate NYS1 region row - remove it
$ws.Rows(19).Delete()

# Update selection to reflect final active cell
$ws.Range("B18").Select()
